$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P5").Value = 40.853215315603983
$ws.Range("Q5").Value = 108.9815951313261
$ws.Range("P6").Value = 40.858750400512648
$ws.Range("Q6").Value = 108.7688084561179
$ws.Range("P7").Value = 40.416924062800391
$ws.Range("Q7").Value = 112.59673286354899
$ws.Range("P8").Value = 41.036055751361737
$ws.Range("Q8").Value = 109.3773030108905
$ws.Range("P9").Value = 41.50213232938161
$ws.Range("Q9").Value = 106.7052402306214
$ws.Range("P10").Value = 40.308361102210831
$ws.Range("Q10").Value = 112.2970916079436
$ws.Range("P11").Value = 41.063115988465228
$ws.Range("Q11").Value = 110.3929212043562
$ws.Range("P12").Value = 40.275948413969878
$ws.Range("Q12").Value = 113.6996540679052
$ws.Range("P13").Value = 40.650972444729263
$ws.Range("Q13").Value = 109.95288917360659
$ws.Range("P14").Value = 40.502161166292851
$ws.Range("Q14").Value = 109.67027546444589
$ws.Range("P15").Value = 40.294458506888823
$ws.Range("Q15").Value = 113.54431133888529
$ws.Range("P16").Value = 41.283590195450167
$ws.Range("Q16").Value = 105.5434913516976
$ws.Range("P17").Value = 41.262894905479023
$ws.Range("Q17").Value = 107.3324727738629
$ws.Range("P18").Value = 40.386244793335479
$ws.Range("Q18").Value = 113.72840486867391
$ws.Range("P19").Value = 41.072569689202183
$ws.Range("Q19").Value = 106.89323510570151
$ws.Range("P20").Value = 41.060696892021788
$ws.Range("Q20").Value = 107.0167392696989
$ws.Range("P21").Value = 40.798438000640807
$ws.Range("Q21").Value = 108.9946828955797
$ws.Range("P22").Value = 40.544508170458187
$ws.Range("Q22").Value = 105.9071748878924
$ws.Range("P23").Value = 40.692947773149619
$ws.Range("Q23").Value = 113.0381550288277
$ws.Range("P24").Value = 40.733973085549501
$ws.Range("Q24").Value = 109.75424087123641
$ws.Range("P25").Value = 41.173715155398909
$ws.Range("Q25").Value = 106.21046764894299
$ws.Range("P26").Value = 40.838711951297661
$ws.Range("Q26").Value = 107.90310698270341
$ws.Range("P27").Value = 40.936094200576733
$ws.Range("Q27").Value = 106.2937411915439
$ws.Range("P28").Value = 41.332866068567768
$ws.Range("Q28").Value = 105.5442921204356
$ws.Range("P29").Value = 40.827835629605893
$ws.Range("Q29").Value = 107.0605637411916
$ws.Range("P30").Value = 40.839778917013767
$ws.Range("Q30").Value = 105.99584240871241
$ws.Range("P31").Value = 40.723522909323933
$ws.Range("Q31").Value = 106.153196668802
$ws.Range("P32").Value = 41.11888176866389
$ws.Range("Q32").Value = 105.5787700192185
$ws.Range("P33").Value = 41.073886574815774
$ws.Range("Q33").Value = 109.36764253683531
$ws.Range("P34").Value = 40.840341236783082
$ws.Range("Q34").Value = 110.2749775784753
$ws.Range("P35").Value = 40.353617430310798
$ws.Range("Q35").Value = 116.54948110185779
$ws.Range("P36").Value = 40.660355655238703
$ws.Range("Q36").Value = 112.5405509288917
$ws.Range("P37").Value = 40.65428548542134
$ws.Range("Q37").Value = 112.14438821268411
$ws.Range("P38").Value = 41.107518423582178
$ws.Range("Q38").Value = 106.7093914157591
$ws.Range("P39").Value = 41.073112784363992
$ws.Range("Q39").Value = 108.79917360666239
$ws.Range("P40").Value = 40.658123998718359
$ws.Range("Q40").Value = 109.6595259449071
$ws.Range("P41").Value = 40.1263889778917
$ws.Range("Q41").Value = 113.15373478539399
$ws.Range("P42").Value = 41.293018263377121
$ws.Range("Q42").Value = 104.747732222934
$ws.Range("P43").Value = 41.497281320089719
$ws.Range("Q43").Value = 104.9949199231262
$ws.Range("P44").Value = 41.45328580583147
$ws.Range("Q44").Value = 104.20614349775779
$ws.Range("P45").Value = 41.141738224927913
$ws.Range("Q45").Value = 105.4260794362588
$ws.Range("P46").Value = 40.202502403075933
$ws.Range("Q46").Value = 112.0481934657271
$ws.Range("P47").Value = 40.088466837552062
$ws.Range("Q47").Value = 117.9945868033312
$ws.Range("P48").Value = 40.652596924062799
$ws.Range("Q48").Value = 112.2269506726457
$ws.Range("P49").Value = 40.406108619032359
$ws.Range("Q49").Value = 112.9717937219731
$ws.Range("P50").Value = 41.150791413008641
$ws.Range("Q50").Value = 107.4082190903267
$ws.Range("P51").Value = 40.924032361422618
$ws.Range("Q51").Value = 106.4392568866111
$ws.Range("P52").Value = 41.073474847805187
$ws.Range("Q52").Value = 104.46516976297239
$ws.Range("P53").Value = 40.657516821531559
$ws.Range("Q53").Value = 112.5238565022422
$ws.Range("P54").Value = 40.820836270426142
$ws.Range("Q54").Value = 110.7881165919282

$ws.Range("T10").Select()
